$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Cells.Item(137, 8).Value = 2797.754
$ws.Cells.Item(137, 9).Value = 2397.878
$ws.Cells.Item(137, 10).Value = 3480.875
$ws.Cells.Item(137, 11).Value = 7193.634
$ws.Cells.Item(137, 12).Value = 10442.625
$ws.Cells.Item(137, 13).Value = -4643.634
$ws.Cells.Item(137, 14).Value = -15542.625
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Cells.Item(63, 8).Value = 2942.7
$ws.Cells.Item(63, 9).Value = 2811.5789
$ws.Cells.Item(63, 10).Value = 3169.182
$ws.Cells.Item(63, 11).Value = 2811.5789
$ws.Cells.Item(63, 12).Value = 3169.182
$ws.Cells.Item(63, 13).Value = -2125.5789
$ws.Cells.Item(63, 14).Value = -4541.182
# Row 66
$ws.Cells.Item(66, 8).Value = 2942.7
$ws.Cells.Item(66, 9).Value = 2811.5789
$ws.Cells.Item(66, 10).Value = 3169.182
$ws.Cells.Item(66, 11).Value = 14057.8945
$ws.Cells.Item(66, 12).Value = 15845.91
$ws.Cells.Item(66, 13).Value = -10625.8945
$ws.Cells.Item(66, 14).Value = -22709.91
# Row 102
$ws.Cells.Item(102, 8).Value = 7002.5
$ws.Cells.Item(102, 9).Value = 6005
$ws.Cells.Item(102, 10).Value = 8000
$ws.Cells.Item(102, 11).Value = 6005
$ws.Cells.Item(102, 12).Value = 8000
$ws.Cells.Item(102, 13).Value = -4383
$ws.Cells.Item(102, 14).Value = -11244
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 5642.9165
$ws.Cells.Item(86, 9).Value = 4651
$ws.Cells.Item(86, 10).Value = 7626.75
$ws.Cells.Item(86, 11).Value = 4651
$ws.Cells.Item(86, 12).Value = 7626.75
$ws.Cells.Item(86, 13).Value = -3528
$ws.Cells.Item(86, 14).Value = -9872.75
# Row 89
$ws.Cells.Item(89, 8).Value = 5642.9165
$ws.Cells.Item(89, 9).Value = 4651
$ws.Cells.Item(89, 10).Value = 7626.75
$ws.Cells.Item(89, 11).Value = 23255
$ws.Cells.Item(89, 12).Value = 38133.75
$ws.Cells.Item(89, 13).Value = -17639
$ws.Cells.Item(89, 14).Value = -49365.75
# Row 94
$ws.Cells.Item(94, 8).Value = 965.85
$ws.Cells.Item(94, 9).Value = 928.7222
$ws.Cells.Item(94, 11).Value = 928.7222
$ws.Cells.Item(94, 13).Value = -477.7222
# Row 99
$ws.Cells.Item(99, 8).Value = 2779069.5
$ws.Cells.Item(99, 9).Value = 4630894
$ws.Cells.Item(99, 11).Value = 4630894
$ws.Cells.Item(99, 13).Value = -4629396
# Row 105
$ws.Cells.Item(105, 8).Value = 2030.25
$ws.Cells.Item(105, 9).Value = 2003.3334
$ws.Cells.Item(105, 10).Value = 2111
$ws.Cells.Item(105, 11).Value = 2003.3334
$ws.Cells.Item(105, 12).Value = 2111
$ws.Cells.Item(105, 13).Value = -256.3334
$ws.Cells.Item(105, 14).Value = -5605
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58, 8).Value = 2155.3333
$ws.Cells.Item(58, 9).Value = 2340
$ws.Cells.Item(58, 10).Value = 1841.4
$ws.Cells.Item(58, 11).Value = 2340
$ws.Cells.Item(58, 12).Value = 1841.4
$ws.Cells.Item(58, 13).Value = -2137
$ws.Cells.Item(58, 14).Value = -2247.4
# Row 62
$ws.Cells.Item(62, 8).Value = 2300
$ws.Cells.Item(62, 10).Value = 2500
$ws.Cells.Item(62, 12).Value = 2500
$ws.Cells.Item(62, 14).Value = -3748
# Row 65
$ws.Cells.Item(65, 8).Value = 2300
$ws.Cells.Item(65, 10).Value = 2500
$ws.Cells.Item(65, 12).Value = 12500
$ws.Cells.Item(65, 14).Value = -18740
# Row 105
$ws.Cells.Item(105, 8).Value = 695.55554
$ws.Cells.Item(105, 9).Value = 635
$ws.Cells.Item(105, 10).Value = 816.6667
$ws.Cells.Item(105, 11).Value = 635
$ws.Cells.Item(105, 12).Value = 816.6667
$ws.Cells.Item(105, 13).Value = 1112
$ws.Cells.Item(105, 14).Value = -4310.6667
# Row 122
$ws.Cells.Item(122, 8).Value = 4624.3447
$ws.Cells.Item(122, 9).Value = 4714.5884
$ws.Cells.Item(122, 10).Value = 4496.5
$ws.Cells.Item(122, 11).Value = 14143.7652
$ws.Cells.Item(122, 12).Value = 13489.5
$ws.Cells.Item(122, 13).Value = -11693.7652
$ws.Cells.Item(122, 14).Value = -18389.5
# Row 132
$ws.Cells.Item(132, 8).Value = 1122.7826
$ws.Cells.Item(132, 9).Value = 893.90247
$ws.Cells.Item(132, 10).Value = 2999.6
$ws.Cells.Item(132, 11).Value = 2681.70741
$ws.Cells.Item(132, 12).Value = 8998.799999999999
$ws.Cells.Item(132, 13).Value = -151.70741
$ws.Cells.Item(132, 14).Value = -14058.8
# Row 136
$ws.Cells.Item(136, 8).Value = 2155.3333
$ws.Cells.Item(136, 9).Value = 2340
$ws.Cells.Item(136, 10).Value = 1841.4
$ws.Cells.Item(136, 11).Value = 7020
$ws.Cells.Item(136, 12).Value = 5524.200000000001
$ws.Cells.Item(136, 13).Value = -4470
$ws.Cells.Item(136, 14).Value = -10624.2
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Cells.Item(86, 8).Value = 433.33334
$ws.Cells.Item(86, 10).Value = 500
$ws.Cells.Item(86, 12).Value = 1500
$ws.Cells.Item(86, 14).Value = -3872
# Row 89
$ws.Cells.Item(89, 8).Value = 433.33334
$ws.Cells.Item(89, 10).Value = 500
$ws.Cells.Item(89, 12).Value = 4500
$ws.Cells.Item(89, 14).Value = -16356
# Row 107
$ws.Cells.Item(107, 8).Value = 322.67856
$ws.Cells.Item(107, 9).Value = 371.25
$ws.Cells.Item(107, 10).Value = 286.25
$ws.Cells.Item(107, 11).Value = 1113.75
$ws.Cells.Item(107, 12).Value = 858.75
$ws.Cells.Item(107, 13).Value = 806.25
$ws.Cells.Item(107, 14).Value = -4698.75
# Row 140
$ws.Cells.Item(140, 8).Value = 2054.5
$ws.Cells.Item(140, 9).Value = 1353.0769
$ws.Cells.Item(140, 10).Value = 3357.1428
$ws.Cells.Item(140, 11).Value = 4059.2307
$ws.Cells.Item(140, 12).Value = 10071.4284
$ws.Cells.Item(140, 13).Value = 1120.7693
$ws.Cells.Item(140, 14).Value = -20431.4284
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 2531.5
$ws.Cells.Item(80, 9).Value = 2499
$ws.Cells.Item(80, 11).Value = 2499
$ws.Cells.Item(80, 13).Value = -1501
# Row 83
$ws.Cells.Item(83, 8).Value = 2531.5
$ws.Cells.Item(83, 9).Value = 2499
$ws.Cells.Item(83, 11).Value = 12495
$ws.Cells.Item(83, 13).Value = -7503
# Row 132
$ws.Cells.Item(132, 8).Value = 1664.3846
$ws.Cells.Item(132, 9).Value = 1103.5358
$ws.Cells.Item(132, 10).Value = 3092
$ws.Cells.Item(132, 11).Value = 3310.6074
$ws.Cells.Item(132, 12).Value = 9276
$ws.Cells.Item(132, 13).Value = -780.6074000000003
$ws.Cells.Item(132, 14).Value = -14336
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Cells.Item(100, 8).Value = 2345.258
$ws.Cells.Item(100, 9).Value = 1989.2222
$ws.Cells.Item(100, 10).Value = 2490.9092
$ws.Cells.Item(100, 11).Value = 1989.2222
$ws.Cells.Item(100, 12).Value = 2490.9092
$ws.Cells.Item(100, 13).Value = -1448.2222
$ws.Cells.Item(100, 14).Value = -3572.9092
# Row 132
$ws.Cells.Item(132, 8).Value = 7922.7036
$ws.Cells.Item(132, 9).Value = 4795.15
$ws.Cells.Item(132, 10).Value = 16858.572
$ws.Cells.Item(132, 11).Value = 14385.45
$ws.Cells.Item(132, 12).Value = 50575.716
$ws.Cells.Item(132, 13).Value = -11855.45
$ws.Cells.Item(132, 14).Value = -55635.716
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 2000
$ws.Cells.Item(62, 9).Value = 2000
$ws.Cells.Item(62, 10).Value = 2000
$ws.Cells.Item(62, 11).Value = 2000
$ws.Cells.Item(62, 12).Value = 2000
$ws.Cells.Item(62, 13).Value = -1376
$ws.Cells.Item(62, 14).Value = -3248
# Row 65
$ws.Cells.Item(65, 8).Value = 2000
$ws.Cells.Item(65, 9).Value = 2000
$ws.Cells.Item(65, 10).Value = 2000
$ws.Cells.Item(65, 11).Value = 10000
$ws.Cells.Item(65, 12).Value = 10000
$ws.Cells.Item(65, 13).Value = -6880
$ws.Cells.Item(65, 14).Value = -16240
# Row 81
$ws.Cells.Item(81, 8).Value = 1390.1818
$ws.Cells.Item(81, 9).Value = 1221.3334
$ws.Cells.Item(81, 10).Value = 2150
$ws.Cells.Item(81, 11).Value = 2442.6668
$ws.Cells.Item(81, 12).Value = 4300
$ws.Cells.Item(81, 13).Value = -1381.6668
$ws.Cells.Item(81, 14).Value = -6422
# Row 84
$ws.Cells.Item(84, 8).Value = 1390.1818
$ws.Cells.Item(84, 9).Value = 1221.3334
$ws.Cells.Item(84, 10).Value = 2150
$ws.Cells.Item(84, 11).Value = 12213.334
$ws.Cells.Item(84, 12).Value = 21500
$ws.Cells.Item(84, 13).Value = -6909.333999999999
$ws.Cells.Item(84, 14).Value = -32108
# Row 96
$ws.Cells.Item(96, 8).Value = 715571.9
$ws.Cells.Item(96, 9).Value = 1500.5
$ws.Cells.Item(96, 10).Value = 5000000
$ws.Cells.Item(96, 11).Value = 1500.5
$ws.Cells.Item(96, 12).Value = 5000000
$ws.Cells.Item(96, 13).Value = -127.5
$ws.Cells.Item(96, 14).Value = -5002746
